$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix the rsId pattern string: drop the surrounding literal quotes ---
$ws.Range("I57").Value = "^rs[0-9]{1,}$"

# Fill color used to highlight "done" rows: Gold, Accent 4, Lighter 60% (#FFE699)
$FillColor = 10086143

# --- Highlight rows + fill in 'nullable' (H) column booleans ---
$ws.Range("H14").Value = $true
$ws.Range("A14:C14").Interior.Color = $FillColor
$ws.Range("H14").Interior.Color = $FillColor
$ws.Range("J14").Interior.Color = $FillColor

$ws.Range("H15").Value = $true
$ws.Range("A15:C15").Interior.Color = $FillColor
$ws.Range("H15").Interior.Color = $FillColor
$ws.Range("J15").Interior.Color = $FillColor

$ws.Range("H17").Value = $false
$ws.Range("A17:B17").Interior.Color = $FillColor
$ws.Range("H17").Interior.Color = $FillColor
$ws.Range("J17").Interior.Color = $FillColor

$ws.Range("H18").Value = $false
$ws.Range("A18:B18").Interior.Color = $FillColor
$ws.Range("H18").Interior.Color = $FillColor
$ws.Range("J18").Interior.Color = $FillColor

$ws.Range("H19").Value = $false
$ws.Range("A19:B19").Interior.Color = $FillColor
$ws.Range("H19").Interior.Color = $FillColor
$ws.Range("J19").Interior.Color = $FillColor

$ws.Range("H25").Value = $false
$ws.Range("A25:B25").Interior.Color = $FillColor
$ws.Range("H25").Interior.Color = $FillColor
$ws.Range("J25").Interior.Color = $FillColor

$ws.Range("H26").Value = $true
$ws.Range("A26:B26").Interior.Color = $FillColor
$ws.Range("H26").Interior.Color = $FillColor
$ws.Range("J26").Interior.Color = $FillColor

$ws.Range("H27").Value = $false
$ws.Range("A27:C27").Interior.Color = $FillColor
$ws.Range("E27:F27").Interior.Color = $FillColor
$ws.Range("H27").Interior.Color = $FillColor
$ws.Range("J27").Interior.Color = $FillColor

$ws.Range("H30").Value = $true
$ws.Range("A30:C30").Interior.Color = $FillColor
$ws.Range("H30").Interior.Color = $FillColor
$ws.Range("J30").Interior.Color = $FillColor

$ws.Range("H33").Value = $true
$ws.Range("A33:C33").Interior.Color = $FillColor
$ws.Range("H33").Interior.Color = $FillColor
$ws.Range("J33").Interior.Color = $FillColor

$ws.Range("H34").Value = $true
$ws.Range("A34:C34").Interior.Color = $FillColor
$ws.Range("H34").Interior.Color = $FillColor
$ws.Range("J34").Interior.Color = $FillColor

$ws.Range("H37").Value = $false
$ws.Range("A37:C37").Interior.Color = $FillColor
$ws.Range("F37").Interior.Color = $FillColor
$ws.Range("H37").Interior.Color = $FillColor
$ws.Range("J37").Interior.Color = $FillColor

$ws.Range("H38").Value = $false
$ws.Range("A38:C38").Interior.Color = $FillColor
$ws.Range("E38").Interior.Color = $FillColor
$ws.Range("H38").Interior.Color = $FillColor
$ws.Range("J38").Interior.Color = $FillColor

$ws.Range("H40").Value = $false
$ws.Range("A40:B40").Interior.Color = $FillColor
$ws.Range("H40").Interior.Color = $FillColor
$ws.Range("J40").Interior.Color = $FillColor

$ws.Range("H43").Value = $false
$ws.Range("A43:C43").Interior.Color = $FillColor
$ws.Range("H43").Interior.Color = $FillColor
$ws.Range("J43").Interior.Color = $FillColor

$ws.Range("H45").Value = $false
$ws.Range("A45:C45").Interior.Color = $FillColor
$ws.Range("E45").Interior.Color = $FillColor
$ws.Range("H45").Interior.Color = $FillColor
$ws.Range("J45").Interior.Color = $FillColor

$ws.Range("H46").Value = $false
$ws.Range("A46:B46").Interior.Color = $FillColor
$ws.Range("H46").Interior.Color = $FillColor
$ws.Range("J46").Interior.Color = $FillColor

$ws.Range("H54").Value = $false
$ws.Range("A54:C54").Interior.Color = $FillColor
$ws.Range("H54").Interior.Color = $FillColor
$ws.Range("J54").Interior.Color = $FillColor

$ws.Range("H56").Value = $false
$ws.Range("A56:C56").Interior.Color = $FillColor
$ws.Range("H56").Interior.Color = $FillColor
$ws.Range("J56").Interior.Color = $FillColor

$ws.Range("H57").Value = $true
$ws.Range("A57:C57").Interior.Color = $FillColor
$ws.Range("H57:J57").Interior.Color = $FillColor

# --- Move the frozen-pane viewport / active selection to where editing left off ---
$ws.Activate()
$ws.Range("I57").Select()
